$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: insert a collapsed "_GoBack" bookmark at the very start of the
# document (before the first run of paragraph 1).
#
# The engine mis-places a bookmark added directly on a zero-length Range(0,0)
# (bookmarkEnd ends up before the next paragraph's run instead of right next
# to bookmarkStart). Work around it: insert a temporary placeholder
# character at position 0, add the bookmark right after it (a non-zero
# position), then delete the placeholder again - the bookmark collapses to
# the correct spot and stays there.
# ---------------------------------------------------------------------------
$r0 = $d.Range(0, 0)
$r0.InsertBefore("X")
$rb = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $rb)
$rx = $d.Range(0, 1)
$rx.Delete()

# ---------------------------------------------------------------------------
# Change 2: split "There were total of 4114 Kickstarter campaigns." into
# three runs: "There were " / "a " / "total of 4114 Kickstarter campaigns."
# while leaving the other runs in the same paragraph (" Overall data..." and
# "close to...") untouched.
#
# A plain Find/Replace or Range.Delete/InsertBefore on this paragraph makes
# the engine re-normalize (merge) every same-formatted run in the paragraph,
# losing the existing run boundaries. Instead, rebuild the whole paragraph
# in one shot via Range.InsertXML, explicitly keeping the untouched runs
# (with their original rsidR) exactly as they were.
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndex $d "There were total of 4114 Kickstarter campaigns."
$paraRange = $d.Paragraphs($idx2).Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7F8E5B99" w14:textId="039607B7" w:rsidR="00A72A9D" w:rsidRDefault="00A72A9D" w:rsidP="00A72A9D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">There were </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:t>total of 4114 Kickstarter campaigns.</w:t></w:r><w:r w:rsidR="006F06BA"><w:t xml:space="preserve"> Overall data of 8 years shows that success percentage is </w:t></w:r><w:r w:rsidR="00EC034C"><w:t>close to the &#8220;Failed&#8221; and &#8220;canceled&#8221; added together.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraRange.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: remove the "_GoBack" bookmark that used to sit in the middle of
# "... failures a|s shown below." and rejoin the surrounding text into a
# single run "... failures as shown below." Rebuild this paragraph the same
# way, via InsertXML, to avoid the engine merging/losing neighbouring runs.
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndex $d "However if we look closely"
$paraRange2 = $d.Paragraphs($idx3).Range
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6D8A0463" w14:textId="1610E698" w:rsidR="00EC034C" w:rsidRDefault="001C4A9B" w:rsidP="00A72A9D"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>However if we look closely, the success and failure is determined by some of the below factors.</w:t></w:r><w:r w:rsidR="004A339A"><w:t xml:space="preserve"> Some categories and sub-categories are close to 100% successful, and some are close to 100% failures as shown below.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraRange2.InsertXML($xml3)

Write-Output "done"
